$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '97.345.77'
$ws.Range("E2").Value = '  +0.49%  '

# Row 3
$ws.Range("D3").Value = '3.700.72'
$ws.Range("E3").Value = '  -0.16%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '2.17'
$ws.Range("E5").Value = '  +14.02%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '236.53'
$ws.Range("E6").Value = '  -1.79%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '656.12'
$ws.Range("E7").Value = '  -0.07%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.440'
$ws.Range("E8").Value = '  +2.71%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.15'
$ws.Range("E9").Value = '  +5.45%  '

# Row 10
$ws.Range("E10").Value = '  -0.09%  '

# Row 11
$ws.Range("D11").Value = '3.698.77'
$ws.Range("E11").Value = '  -0.10%  '

# Row 12
$ws.Range("B12").Value = 'Avalanche'
$ws.Range("C12").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.05'
$ws.Range("E12").Value = '  -0.81%  '

# Row 13
$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000312'
$ws.Range("E13").Value = '  +15.24%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.208'
$ws.Range("E14").Value = '  +0.29%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.84'
$ws.Range("E15").Value = '  -1.03%  '

# Row 16
$ws.Range("D16").Value = '4.389.64'
$ws.Range("E16").Value = '  -0.19%  '

# Row 17
$ws.Range("D17").Value = '97.005.76'
$ws.Range("E17").Value = '  +0.26%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.87'
$ws.Range("E18").Value = '  -2.09%  '

# Row 19
$ws.Range("D19").Value = '3.700.92'
$ws.Range("E19").Value = '  -0.03%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.14'
$ws.Range("E20").Value = '  +1.74%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.92'
$ws.Range("E21").Value = '  -1.60%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.543'
$ws.Range("E22").Value = '  +2.61%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '526.12'
$ws.Range("E23").Value = '  -0.09%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.44'
$ws.Range("E24").Value = '  -2.22%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000222'
$ws.Range("E25").Value = '  +8.29%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '117.51'
$ws.Range("E26").Value = '  +14.73%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.91'
$ws.Range("E27").Value = '  -2.61%  '

# Row 28
$ws.Range("E28").Value = '  +21.46%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '13.43'
$ws.Range("E29").Value = '  -0.24%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.69'
$ws.Range("E30").Value = '  +0.06%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.03'
$ws.Range("E31").Value = '  -1.42%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.998'
$ws.Range("E32").Value = '  -0.18%  '

# Row 33
$ws.Range("E33").Value = '  +1.39%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '33.05'
$ws.Range("E34").Value = '  +0.46%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.81'
$ws.Range("E35").Value = '  -4.04%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.997'
$ws.Range("E36").Value = '  -0.54%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.597'
$ws.Range("E37").Value = '  -1.30%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '633.66'
$ws.Range("E38").Value = '  -3.73%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.74'
$ws.Range("E39").Value = '  -2.90%  '

# Row 41
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.507'
$ws.Range("E41").Value = '  +14.63%  '

# Row 42
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.166'
$ws.Range("E42").Value = '  +1.76%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.86'
$ws.Range("E43").Value = '  -4.55%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.27'

# Row 45
$ws.Range("E45").Value = '  +0.11%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.961'
$ws.Range("E46").Value = '  -1.34%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0455'
$ws.Range("E47").Value = '  -1.14%  '

# Row 48
$ws.Range("E48").Value = '  +1.56%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.81'
$ws.Range("E49").Value = '  +0.70%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.63'
$ws.Range("E50").Value = '  +0.03%  '

# Row 51
$ws.Range("E51").Value = '  +2.87%  '
